$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the Hortaliza / Alcachofa price records (Terminal La Palmera de La Serena).
# Each data row (2-19) is updated in place with this weeks values for Fecha, Variedad,
# Volumen, Precio minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg and Kg o Unidades.

# Row 2
$ws.Range("D2").Value = 44484
$ws.Range("H2").Value = 'Española'
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9500
$ws.Range("N2").Value = '$/caja 30 unidades'
$ws.Range("O2").Value = 'Provincia del Elquí'
$ws.Range("P2").Value = 317
$ws.Range("Q2").Value = 30

# Row 3
$ws.Range("D3").Value = 44858
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 9500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9750
$ws.Range("P3").Value = 325

# Row 4
$ws.Range("D4").Value = 44420
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("N4").Value = '$/caja 40 unidades'
$ws.Range("O4").Value = 'Provincia de Limarí'
$ws.Range("P4").Value = 362
$ws.Range("Q4").Value = 40

# Row 5
$ws.Range("D5").Value = 44420
$ws.Range("H5").Value = 'Madrigal'
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 13500
$ws.Range("N5").Value = '$/caja 40 unidades'
$ws.Range("P5").Value = 338
$ws.Range("Q5").Value = 40

# Row 6
$ws.Range("D6").Value = 44729
$ws.Range("H6").Value = 'Madrigal'
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 16500
$ws.Range("N6").Value = '$/caja 40 unidades'
$ws.Range("O6").Value = 'Provincia del Elquí'
$ws.Range("P6").Value = 412
$ws.Range("Q6").Value = 40

# Row 7
$ws.Range("H7").Value = 'Española'
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14500
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("O7").Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 483
$ws.Range("Q7").Value = 30

# Row 8
$ws.Range("D8").Value = 44790
$ws.Range("H8").Value = 'Madrigal'
$ws.Range("K8").Value = 11500
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 11750
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("P8").Value = 294
$ws.Range("Q8").Value = 40

# Row 9
$ws.Range("D9").Value = 44701
$ws.Range("H9").Value = 'Española'
$ws.Range("K9").Value = 19000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19500
$ws.Range("N9").Value = '$/caja 30 unidades'
$ws.Range("P9").Value = 650
$ws.Range("Q9").Value = 30

# Row 10
$ws.Range("D10").Value = 44784
$ws.Range("H10").Value = 'Madrigal'
$ws.Range("J10").Value = 520
$ws.Range("K10").Value = 11500
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11750
$ws.Range("N10").Value = '$/caja 40 unidades'
$ws.Range("P10").Value = 294
$ws.Range("Q10").Value = 40

# Row 11
$ws.Range("D11").Value = 44855
$ws.Range("J11").Value = 540
$ws.Range("K11").Value = 9500
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9750
$ws.Range("O11").Value = 'Provincia del Elquí'
$ws.Range("P11").Value = 325

# Row 12
$ws.Range("D12").Value = 44687
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 19000
$ws.Range("M12").Value = 18500
$ws.Range("P12").Value = 617

# Row 13
$ws.Range("D13").Value = 44427
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12500
$ws.Range("P13").Value = 312

# Row 14
$ws.Range("D14").Value = 44426
$ws.Range("H14").Value = 'Española'
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 11500
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11750
$ws.Range("N14").Value = '$/caja 30 unidades'
$ws.Range("O14").Value = 'Provincia de Limarí'
$ws.Range("P14").Value = 392
$ws.Range("Q14").Value = 30

# Row 15
$ws.Range("H15").Value = 'Madrigal'
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 12500
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12750
$ws.Range("N15").Value = '$/caja 40 unidades'
$ws.Range("P15").Value = 319
$ws.Range("Q15").Value = 40

# Row 16
$ws.Range("D16").Value = 44498
$ws.Range("H16").Value = 'Española'
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 8500
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = 8750
$ws.Range("N16").Value = '$/caja 30 unidades'
$ws.Range("P16").Value = 292
$ws.Range("Q16").Value = 30

# Row 17
$ws.Range("D17").Value = 44767
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("O17").Value = 'Provincia de Limarí'
$ws.Range("P17").Value = 362

# Row 18
$ws.Range("D18").Value = 44839
$ws.Range("H18").Value = 'Española'
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 13000
$ws.Range("M18").Value = 12500
$ws.Range("N18").Value = '$/caja 30 unidades'
$ws.Range("O18").Value = 'Provincia del Elquí'
$ws.Range("P18").Value = 417
$ws.Range("Q18").Value = 30

# Row 19
$ws.Range("D19").Value = 44438
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 11000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 11500
$ws.Range("P19").Value = 383

